# Auto-generated edit script: updates live market price columns (H-N)
# across multiple worksheets, per scheduled Chocobo Profits data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 188.36842
$ws.Range("I33").Value = 104
$ws.Range("J33").Value = 504.75
$ws.Range("K33").Value = 104
$ws.Range("L33").Value = 504.75
$ws.Range("M33").Value = 125
$ws.Range("N33").Value = -962.75

$ws.Range("H41").Value = 667.5625
$ws.Range("I41").Value = 440.2
$ws.Range("J41").Value = 770.9091
$ws.Range("K41").Value = 440.2
$ws.Range("L41").Value = 770.9091
$ws.Range("M41").Value = -0.1999999999999886
$ws.Range("N41").Value = -1650.9091

$ws.Range("H53").Value = 462.22223
$ws.Range("I53").Value = 357.58334
$ws.Range("K53").Value = 357.58334
$ws.Range("M53").Value = 279.41666

$ws.Range("H135").Value = 394.86957
$ws.Range("I135").Value = 365.72726
$ws.Range("K135").Value = 3291.54534
$ws.Range("M135").Value = -756.5453400000001

$ws.Range("H137").Value = 1930.0416
$ws.Range("I137").Value = 1153.8235
$ws.Range("J137").Value = 3815.1428
$ws.Range("K137").Value = 3461.4705
$ws.Range("L137").Value = 11445.4284
$ws.Range("M137").Value = -911.4704999999999
$ws.Range("N137").Value = -16545.4284

$ws.Range("H138").Value = 2773.26
$ws.Range("I138").Value = 959.6667
$ws.Range("J138").Value = 3255.3545
$ws.Range("K138").Value = 2879.0001
$ws.Range("L138").Value = 9766.0635
$ws.Range("M138").Value = 2260.9999
$ws.Range("N138").Value = -20046.0635

$ws.Range("H141").Value = 5139.125
$ws.Range("I141").Value = 5378.7607
$ws.Range("J141").Value = 4036.8
$ws.Range("K141").Value = 16136.2821
$ws.Range("L141").Value = 12110.4
$ws.Range("M141").Value = -10956.2821
$ws.Range("N141").Value = -22470.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3051.0466
$ws.Range("I32").Value = 2726.4487
$ws.Range("K32").Value = 2726.4487
$ws.Range("M32").Value = -2439.4487

$ws.Range("H74").Value = 2184.7258
$ws.Range("I74").Value = 2020.5636
$ws.Range("K74").Value = 2020.5636
$ws.Range("M74").Value = -1146.5636

$ws.Range("H77").Value = 2184.7258
$ws.Range("I77").Value = 2020.5636
$ws.Range("K77").Value = 10102.818
$ws.Range("M77").Value = -5734.817999999999

$ws.Range("H122").Value = 1581.06
$ws.Range("I122").Value = 1135.3077
$ws.Range("J122").Value = 3161.4546
$ws.Range("K122").Value = 3405.9231
$ws.Range("L122").Value = 9484.363799999999
$ws.Range("M122").Value = -955.9231
$ws.Range("N122").Value = -14384.3638

$ws.Range("H134").Value = 39998
$ws.Range("J134").Value = 39998
$ws.Range("L134").Value = 39998
$ws.Range("N134").Value = -50138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1610
$ws.Range("I134").Value = 921.21313
$ws.Range("J134").Value = 3821.3684
$ws.Range("K134").Value = 2763.63939
$ws.Range("L134").Value = 11464.1052
$ws.Range("M134").Value = -228.6393899999998
$ws.Range("N134").Value = -16534.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9436621
$ws.Range("I31").Value = 1392.8846
$ws.Range("J31").Value = 18522398
$ws.Range("K31").Value = 1392.8846
$ws.Range("L31").Value = 18522398
$ws.Range("M31").Value = -1097.8846
$ws.Range("N31").Value = -18522988

$ws.Range("H34").Value = 9436621
$ws.Range("I34").Value = 1392.8846
$ws.Range("J34").Value = 18522398
$ws.Range("K34").Value = 1392.8846
$ws.Range("L34").Value = 18522398
$ws.Range("M34").Value = -1190.8846
$ws.Range("N34").Value = -18522802

$ws.Range("H58").Value = 1602.573
$ws.Range("I58").Value = 1361.1
$ws.Range("K58").Value = 1361.1
$ws.Range("M58").Value = -1158.1

$ws.Range("H132").Value = 3027.6206
$ws.Range("I132").Value = 2518.5557
$ws.Range("J132").Value = 9900
$ws.Range("K132").Value = 7555.6671
$ws.Range("L132").Value = 29700
$ws.Range("M132").Value = -5025.6671
$ws.Range("N132").Value = -34760

$ws.Range("H136").Value = 1602.573
$ws.Range("I136").Value = 1361.1
$ws.Range("K136").Value = 4083.3
$ws.Range("M136").Value = -1533.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1301.0278
$ws.Range("I5").Value = 291.14285
$ws.Range("J5").Value = 2714.8667
$ws.Range("K5").Value = 873.4285500000001
$ws.Range("L5").Value = 8144.6001
$ws.Range("M5").Value = -761.4285500000001
$ws.Range("N5").Value = -8368.6001

$ws.Range("H121").Value = 3148.4146
$ws.Range("I121").Value = 215
$ws.Range("J121").Value = 3298.8462
$ws.Range("K121").Value = 645
$ws.Range("L121").Value = 9896.5386
$ws.Range("M121").Value = 665
$ws.Range("N121").Value = -12516.5386

$ws.Range("H131").Value = 5435619
$ws.Range("J131").Value = 859.6
$ws.Range("L131").Value = 2578.8
$ws.Range("N131").Value = -12658.8

$ws.Range("H135").Value = 1301.0278
$ws.Range("I135").Value = 291.14285
$ws.Range("J135").Value = 2714.8667
$ws.Range("K135").Value = 2620.28565
$ws.Range("L135").Value = 24433.8003
$ws.Range("M135").Value = -85.28565000000026
$ws.Range("N135").Value = -29503.8003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2418.8484
$ws.Range("I122").Value = 1245.7273
$ws.Range("K122").Value = 3737.1819
$ws.Range("M122").Value = -1287.1819

$ws.Range("H132").Value = 2268.2092
$ws.Range("I132").Value = 1467.4242
$ws.Range("J132").Value = 4910.8
$ws.Range("K132").Value = 4402.2726
$ws.Range("L132").Value = 14732.4
$ws.Range("M132").Value = -1872.2726
$ws.Range("N132").Value = -19792.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2076.7256
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2076.7256
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2076.7256
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2300.7256

$ws.Range("H46").Value = 1336.909
$ws.Range("I46").Value = 1239.091
$ws.Range("J46").Value = 1434.7273
$ws.Range("K46").Value = 1239.091
$ws.Range("L46").Value = 1434.7273
$ws.Range("M46").Value = -1051.091
$ws.Range("N46").Value = -1810.7273

$ws.Range("H75").Value = 38000
$ws.Range("I75").Value = 2000
$ws.Range("K75").Value = 2000
$ws.Range("M75").Value = -1064

$ws.Range("H78").Value = 38000
$ws.Range("I78").Value = 2000
$ws.Range("K78").Value = 6000
$ws.Range("M78").Value = -1320

$ws.Range("H80").Value = 47437.5
$ws.Range("J80").Value = 47437.5
$ws.Range("L80").Value = 47437.5
$ws.Range("N80").Value = -49683.5

$ws.Range("H83").Value = 47437.5
$ws.Range("J83").Value = 47437.5
$ws.Range("L83").Value = 142312.5
$ws.Range("N83").Value = -153544.5

$ws.Range("H132").Value = 10557.55
$ws.Range("I132").Value = 11830.5
$ws.Range("J132").Value = 8648.125
$ws.Range("K132").Value = 35491.5
$ws.Range("L132").Value = 25944.375
$ws.Range("M132").Value = -32961.5
$ws.Range("N132").Value = -31004.375

$ws.Range("H136").Value = 3431.16
$ws.Range("I136").Value = 1840.7368
$ws.Range("J136").Value = 8467.5
$ws.Range("K136").Value = 5522.2104
$ws.Range("L136").Value = 25402.5
$ws.Range("M136").Value = -2972.2104
$ws.Range("N136").Value = -30502.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2937.2903
$ws.Range("I122").Value = 1967.1904
$ws.Range("J122").Value = 4974.5
$ws.Range("K122").Value = 5901.5712
$ws.Range("L122").Value = 14923.5
$ws.Range("M122").Value = -3451.5712
$ws.Range("N122").Value = -19823.5

$ws.Range("H136").Value = 1785.2307
$ws.Range("I136").Value = 527.9792
$ws.Range("K136").Value = 1583.9376
$ws.Range("M136").Value = 966.0624
